$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "reviews_count" column (E) is removed entirely; everything to its
# right (reviews_average .. latest_review_date, originally F:K) shifts
# one column to the left (becoming E:J).
$ws.Columns("E").Delete()
